$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$eValues = @(0.9926147277508489,0.9936372048519304,0.9942998659930995,0.9945786998346017,0.9946255319796338,0.9943035907982488,0.9929600610674301,0.9906006454969559,0.989033133672735,0.988355674866747,0.9881042295826724,0.9881581567098651,0.9883348863814464,0.9884438009545853,0.9890781214508737,0.989476357848556,0.9897087662937556,0.9897880325774034,0.9894336180360679,0.9882828385668253,0.9875604150241495,0.9879432794643023,0.9894529299347244,0.9912096547607049)
$lValues = @(0.9955398523336033,0.9963617723202692,0.9968940712668345,0.997117960005301,0.9971555583673453,0.9968970624462087,0.995817528259106,0.9939188001724441,0.9926553831429383,0.9921088820399291,0.9919059725120875,0.9919494934313052,0.9920921077337197,0.9921799884222134,0.9926916645766087,0.9930127773699352,0.9932001317071769,0.9932640239640975,0.9929783193494216,0.9920501090198105,0.9914670000341481,0.991776070289318,0.9929938892766442,0.9944092447426414)

for ($i = 0; $i -lt $eValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $eValues[$i]
    $ws.Cells.Item($row, 12).Value = $lValues[$i]
}
